$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two changed point values in row 2 (C2, D2)
$ws.Range("C2").Value = "Point_6"
$ws.Range("D2").Value = "Point_14"

# Move the active selection from D10 to D11
$ws.Range("D11").Select()
